# Fruta / hortaliza, semanal
# Insert a new weekly record for Berenjena (row 95) and push the
# existing records (previously rows 95-118) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 95:118 down to 96:119, leaving a blank row 95 for the new record.
$ws.Rows("95").Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A95").Value = 6
$ws.Range("B95").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C95").Value = "Metropolitana"
$ws.Range("D95").Value = 44508
$ws.Range("E95").Value = 13
$ws.Range("F95").Value = 100112001
$ws.Range("G95").Value = "Berenjena"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 400
$ws.Range("K95").Value = 7000
$ws.Range("L95").Value = 8000
$ws.Range("M95").Value = 7425
$ws.Range("N95").Value = "$/caja 50 unidades"
$ws.Range("O95").Value = "Región de Arica y Parinacota"
$ws.Range("P95").Value = 148
$ws.Range("Q95").Value = 50
$ws.Range("R95").Value = "Hortaliza"
